# Applies the recorded edit to the "Artfynd" sheet:
#  - Rows 31-43: the species/observation data in columns A,B,D,E,F,G,H,I,J is
#    permuted between rows (the underlying records were re-sorted upstream),
#    and columns Q/R (Ost/Nord easting-northing coordinates) are rewritten
#    as rounded integers, following the same row mapping.
#  - Rows 9 and 30: columns Q/R are simply rounded to whole numbers in place.
#  - Rows 9 and 30-43: the Starttid/Sluttid columns (Z, AB), which only ever
#    contained the placeholder value "00:00", are cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row -> source-row mapping for the species/observation columns (A,B,D,E,F,G,H,I,J)
# and the Q/R coordinate columns, for rows 31-43.
$mapping = @{
    31 = 39
    32 = 37
    33 = 33
    34 = 31
    35 = 34
    36 = 32
    37 = 43
    38 = 36
    39 = 41
    40 = 35
    41 = 40
    42 = 42
    43 = 38
}

$cols = @("A","B","D","E","F","G","H","I","J")

# --- Step 1: snapshot all current values for rows 31-43 before overwriting anything ---
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $rowData["Q"] = $ws.Range("Q$row").Value2
    $rowData["R"] = $ws.Range("R$row").Value2
    $snapshot[$row] = $rowData
}

# --- Step 2: write the permuted values into rows 31-43 ---
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $src = $snapshot[$srcRow]

    foreach ($col in $cols) {
        $val = $src[$col]
        if ($null -eq $val) {
            $ws.Range("$col$destRow").Value = ""
        } else {
            $ws.Range("$col$destRow").Value = $val
        }
    }

    $ws.Range("Q$destRow").Value = [math]::Round([double]$src["Q"])
    $ws.Range("R$destRow").Value = [math]::Round([double]$src["R"])
}

# --- Step 3: rows 9 and 30 just get Q/R rounded in place ---
foreach ($row in @(9, 30)) {
    $q = $ws.Range("Q$row").Value2
    $r = $ws.Range("R$row").Value2
    $ws.Range("Q$row").Value = [math]::Round([double]$q)
    $ws.Range("R$row").Value = [math]::Round([double]$r)
}

# --- Step 4: clear the Starttid/Sluttid ("00:00") cells for row 9 and rows 30-43 ---
$timeRows = @(9) + @(30..43)
foreach ($row in $timeRows) {
    $ws.Range("Z$row").ClearContents()
    $ws.Range("AB$row").ClearContents()
}
